# Excel Change Because of Nainsi
#
# The "Contract" sheet had its contract id (cell A3) updated from the old
# cancelled contract number "CS136744A9" to the new one "CE002662A9".
# The Contract sheet becomes the active sheet/tab, with the cell below the
# edited one (A4) left selected, matching how Excel leaves the selection
# after typing a value into a cell and pressing Enter.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Contract")
$ws.Activate()

$ws.Range("A3").Value = "CE002662A9"

$ws.Range("A4").Select()
